$wb = $excel.ActiveWorkbook

# Sheet "20_Properties of Circles": A1 held the full sheet title as text;
# replace it with just the leading question-set number (20).
$ws1 = $wb.Worksheets.Item("20_Properties of Circles")
$ws1.Range("A1").Value = 20

# Sheet "21_Prop of Tangent to Circle": same treatment, number 21.
$ws2 = $wb.Worksheets.Item("21_Prop of Tangent to Circle")
$ws2.Range("A1").Value = 21

# Leave sheet 2 with A2 selected, then make sheet 1 the active tab with A2
# selected there too (matches the saved view state in the target workbook).
$ws2.Activate()
$ws2.Range("A2").Select()

$ws1.Activate()
$ws1.Range("A2").Select()
